# daily auto push: 2026-01-18 22:34 UTC
#
# A new sample row needs to be inserted right before the existing row 667
# (date 2026/01/19, day 月, time 4, ranking 201). Every row from the old
# 667 downward shifts down by one (old 667 -> new 668, ..., old 708 ->
# new 709), and the sheet's used-range dimension grows from D708 to D709.
#
# Copy-then-insert row 666 (which already holds 2026/01/19 / 月 as plain
# text, matching the column's existing inline-string-like data) so the
# newly inserted row inherits the same cell formatting/type as its
# neighbours instead of Excel's "looks like a date" auto-conversion that
# a plain string assignment into a blank cell would trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(666).Copy()
$ws.Rows.Item(667).Insert()

# Row 666 already holds "2026/01/19" / "月" / 201 (D), identical to what the
# new row needs, so the copy leaves A667/B667/D667 correct as-is. Only the
# time column differs, so only C667 needs an explicit write. (Re-assigning
# A667/B667 here would make Excel re-interpret the date-like string as a
# real date/number and defeat the copied text formatting.)
$ws.Range("C667").Value = 4
